$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs Thpo Mpl ECs)
$ws.Range("G2").Value = 0.09370333333333335
$ws.Range("H2").Value = 0.28111
$ws.Range("I2").Value = 0.3834601030162819
$ws.Range("J2").Value = 0.3834601030162819
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.033008
$ws.Range("N2").Value = 0.099024
$ws.Range("O2").Value = 0.1900718258980157
$ws.Range("P2").Value = 0.1900718258980157
$ws.Range("Q2").Value = 0.003092959626666668
$ws.Range("R2").Value = 0.02783663664
$ws.Range("S2").Value = 0.07288496193934588
$ws.Range("T2").Value = 0.07288496193934588

# Row 3 (ECs Thpo Mpl FAPs)
$ws.Range("G3").Value = 0.09370333333333335
$ws.Range("H3").Value = 0.28111
$ws.Range("I3").Value = 0.3834601030162819
$ws.Range("J3").Value = 0.3834601030162819
$ws.Range("O3").Value = 0.8099281741019843
$ws.Range("P3").Value = 0.8099281741019843
$ws.Range("Q3").Value = 0.01317962370888889
$ws.Range("R3").Value = 0.11861661338
$ws.Range("S3").Value = 0.310575141076936
$ws.Range("T3").Value = 0.310575141076936

# Row 4 (FAPs Thpo Mpl ECs)
$ws.Range("I4").Value = 0.2360439674363787
$ws.Range("J4").Value = 0.2360439674363787
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.033008
$ws.Range("N4").Value = 0.099024
$ws.Range("O4").Value = 0.1900718258980157
$ws.Range("P4").Value = 0.1900718258980157
$ws.Range("Q4").Value = 0.001903912442666667
$ws.Range("R4").Value = 0.017135211984
$ws.Range("S4").Value = 0.04486530788284425
$ws.Range("T4").Value = 0.04486530788284426

# Row 5 (FAPs Thpo Mpl FAPs)
$ws.Range("I5").Value = 0.2360439674363787
$ws.Range("J5").Value = 0.2360439674363787
$ws.Range("O5").Value = 0.8099281741019843
$ws.Range("P5").Value = 0.8099281741019843
$ws.Range("S5").Value = 0.1911786595535344
$ws.Range("T5").Value = 0.1911786595535345

# Row 6 (MuSCs Thpo Mpl ECs)
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.09297899999999999
$ws.Range("H6").Value = 0.278937
$ws.Range("I6").Value = 0.3804959295473394
$ws.Range("J6").Value = 0.3804959295473395
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.033008
$ws.Range("N6").Value = 0.099024
$ws.Range("O6").Value = 0.1900718258980157
$ws.Range("P6").Value = 0.1900718258980157
$ws.Range("Q6").Value = 0.003069050832
$ws.Range("R6").Value = 0.027621457488
$ws.Range("S6").Value = 0.07232155607582554
$ws.Range("T6").Value = 0.07232155607582554

# Row 7 (MuSCs Thpo Mpl FAPs)
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.09297899999999999
$ws.Range("H7").Value = 0.278937
$ws.Range("I7").Value = 0.3804959295473394
$ws.Range("J7").Value = 0.3804959295473395
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.033008
$ws.Range("N7").Value = 0.099024
$ws.Range("O7").Value = 0.1900718258980157
$ws.Range("P7").Value = 0.1900718258980157
$ws.Range("Q7").Value = 0.013077744294
$ws.Range("R7").Value = 0.117699698646
$ws.Range("S7").Value = 0.3081743734715138
$ws.Range("T7").Value = 0.3081743734715139
